# Update cryptocurrency price/volume data in the active worksheet
# (values scraped nightly by GitHub Actions workflow).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Several "Price" cells contain plain numeric-looking text (e.g. "528.38").
# Force those specific cells to Text format first so Excel keeps them as
# literal strings instead of auto-converting them to numbers.
$textCells = @(
    "D5", "D6", "D7", "D9", "D10", "D11", "D12", "D14", "D16", "D18",
    "D19", "D21", "D22", "D24", "D25", "D26", "D27", "D28", "D29", "D30",
    "D32", "D33", "D34", "D35", "D36", "D38", "D41", "D43", "D44", "D45",
    "D47", "D48", "D49"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values cell by cell.
$ws.Range("D2").Value = "69.636.54"
$ws.Range("E2").Value = "  +1.79%  "
$ws.Range("D3").Value = "3.924.84"
$ws.Range("E3").Value = "  +0.04%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "528.38"
$ws.Range("E5").Value = "  +8.57%  "
$ws.Range("D6").Value = "145.09"
$ws.Range("E6").Value = "  -0.81%  "
$ws.Range("D7").Value = "0.616"
$ws.Range("E7").Value = "  -1.24%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("D9").Value = "0.730"
$ws.Range("E9").Value = "  -0.53%  "
$ws.Range("D10").Value = "0.173"
$ws.Range("E10").Value = "  +3.51%  "
$ws.Range("D11").Value = "0.0000335"
$ws.Range("E11").Value = "  -3.55%  "
$ws.Range("D12").Value = "42.56"
$ws.Range("E12").Value = "  -1.48%  "
$ws.Range("D13").Value = "4.549.77"
$ws.Range("E13").Value = "  +0.11%  "
$ws.Range("D14").Value = "10.31"
$ws.Range("E14").Value = "  -4.17%  "
$ws.Range("D15").Value = "3.933.69"
$ws.Range("E15").Value = "  -0.03%  "
$ws.Range("D16").Value = "1.24"
$ws.Range("E16").Value = "  +8.81%  "
$ws.Range("E17").Value = "  -0.27%  "
$ws.Range("D18").Value = "13.93"
$ws.Range("E18").Value = "  -2.99%  "
$ws.Range("D19").Value = "19.97"
$ws.Range("D20").Value = "69.507.72"
$ws.Range("E20").Value = "  +1.58%  "
$ws.Range("D21").Value = "438.32"
$ws.Range("E21").Value = "  +1.16%  "
$ws.Range("D22").Value = "3.36"
$ws.Range("E22").Value = "  -3.89%  "
$ws.Range("E23").Value = "  -5.69%  "
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").Value = "88.13"
$ws.Range("E24").Value = "  -0.26%  "
$ws.Range("B25").Value = "PancakeSwap"
$ws.Range("C25").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D25").Value = "4.08"
$ws.Range("E25").Value = "  +11.55%  "
$ws.Range("D26").Value = "11.57"
$ws.Range("E26").Value = "  -0.22%  "
$ws.Range("D27").Value = "10.73"
$ws.Range("E27").Value = "  -4.89%  "
$ws.Range("D28").Value = "36.45"
$ws.Range("E28").Value = "  -4.13%  "
$ws.Range("D29").Value = "697.52"
$ws.Range("E29").Value = "  -2.47%  "
$ws.Range("D30").Value = "13.26"
$ws.Range("E30").Value = "  -4.14%  "
$ws.Range("E31").Value = "  -2.59%  "
$ws.Range("D32").Value = "2.84"
$ws.Range("E32").Value = "  -3.36%  "
$ws.Range("D33").Value = "68.04"
$ws.Range("E33").Value = "  +11.58%  "
$ws.Range("D34").Value = "0.448"
$ws.Range("E34").Value = "  +12.73%  "
$ws.Range("B35").Value = "InjectiveProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D35").Value = "40.65"
$ws.Range("E35").Value = "  -2.06%  "
$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D36").Value = "5.96"
$ws.Range("E36").Value = "  -3.15%  "
$ws.Range("D37").Value = "0.0₃0840"
$ws.Range("E37").Value = "  -5.29%  "
$ws.Range("D38").Value = "0.151"
$ws.Range("E38").Value = "  +3.59%  "
$ws.Range("E39").Value = "  -0.06%  "
$ws.Range("E40").Value = "  -0.06%  "
$ws.Range("D41").Value = "0.0484"
$ws.Range("E41").Value = "  -2.26%  "
$ws.Range("E42").Value = "  +3.40%  "
$ws.Range("D43").Value = "2.77"
$ws.Range("E43").Value = "  -8.18%  "
$ws.Range("D44").Value = "2.96"
$ws.Range("E44").Value = "  -5.15%  "
$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D45").Value = "3.38"
$ws.Range("E45").Value = "  -0.32%  "
$ws.Range("E46").Value = "  -0.59%  "
$ws.Range("B47").Value = "Stacks"
$ws.Range("C47").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D47").Value = "3.11"
$ws.Range("E47").Value = "  +10.08%  "
$ws.Range("D48").Value = "3.32"
$ws.Range("E48").Value = "  -2.93%  "
$ws.Range("D49").Value = "144.77"
$ws.Range("D50").Value = "0.0₆0341"
$ws.Range("E50").Value = "  -1.67%  "
$ws.Range("E51").Value = "  -3.49%  "
